$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 7191.7617
$ws.Range("J17").Value = 7456.4
$ws.Range("L17").Value = 22369.2
$ws.Range("N17").Value = -22705.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1250
$ws.Range("I98").Value = 1250
$ws.Range("K98").Value = 1250
$ws.Range("M98").Value = 248

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1959.5454
$ws.Range("I100").Value = 1551.25
$ws.Range("J100").Value = 2192.8572
$ws.Range("K100").Value = 1551.25
$ws.Range("L100").Value = 2192.8572
$ws.Range("M100").Value = -1010.25
$ws.Range("N100").Value = -3274.8572

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H109").Value = 27000
$ws.Range("J109").Value = 27000
$ws.Range("L109").Value = 27000
$ws.Range("N109").Value = -29774

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 321.5
$ws.Range("I115").Value = 321.5
$ws.Range("K115").Value = 964.5
$ws.Range("M115").Value = 602.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 469.85715
$ws.Range("I118").Value = 469.85715
$ws.Range("K118").Value = 1409.57145
$ws.Range("M118").Value = 247.4285500000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 7559.706
$ws.Range("I121").Value = 1215
$ws.Range("J121").Value = 8919.286
$ws.Range("K121").Value = 3645
$ws.Range("L121").Value = 26757.858
$ws.Range("M121").Value = -1898
$ws.Range("N121").Value = -30251.858

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1250
$ws.Range("I122").Value = 1250
$ws.Range("K122").Value = 3750
$ws.Range("M122").Value = -1300

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 29995
$ws.Range("J123").Value = 29995
$ws.Range("L123").Value = 29995
$ws.Range("N123").Value = -39795

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 765.35596
$ws.Range("J129").Value = 793.8
$ws.Range("L129").Value = 2381.4
$ws.Range("N129").Value = -12381.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 15155084
$ws.Range("I135").Value = 624.26666
$ws.Range("K135").Value = 5618.39994
$ws.Range("M135").Value = -3083.39994

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2270.2896
$ws.Range("I138").Value = 1175.6428
$ws.Range("J138").Value = 2908.8333
$ws.Range("K138").Value = 3526.9284
$ws.Range("L138").Value = 8726.499899999999
$ws.Range("M138").Value = 1613.0716
$ws.Range("N138").Value = -19006.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7063.241
$ws.Range("I32").Value = 5694.0146
$ws.Range("K32").Value = 5694.0146
$ws.Range("M32").Value = -5407.0146

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 29413268
$ws.Range("I74").Value = 37037692
$ws.Range("K74").Value = 37037692
$ws.Range("M74").Value = -37036818

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 29413268
$ws.Range("I77").Value = 37037692
$ws.Range("K77").Value = 185188460
$ws.Range("M77").Value = -185184092

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 35354
$ws.Range("J87").Value = 35354
$ws.Range("L87").Value = 35354
$ws.Range("N87").Value = -37850

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H90").Value = 35354
$ws.Range("J90").Value = 35354
$ws.Range("L90").Value = 106062
$ws.Range("N90").Value = -118542

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 863967.25
$ws.Range("I105").Value = 1479.5834
$ws.Range("J105").Value = 1472782
$ws.Range("K105").Value = 1479.5834
$ws.Range("L105").Value = 1472782
$ws.Range("M105").Value = 267.4166
$ws.Range("N105").Value = -1476276

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3847.6667
$ws.Range("I134").Value = 4101.967
$ws.Range("J134").Value = 1304.6666
$ws.Range("K134").Value = 12305.901
$ws.Range("L134").Value = 3913.9998
$ws.Range("M134").Value = -9770.900999999998
$ws.Range("N134").Value = -8983.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 170.4
$ws.Range("I22").Value = 148.42857
$ws.Range("J22").Value = 221.66667
$ws.Range("K22").Value = 148.42857
$ws.Range("L22").Value = 221.66667
$ws.Range("M22").Value = 201.57143
$ws.Range("N22").Value = -921.6666700000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4713.3096
$ws.Range("I31").Value = 2566.7334
$ws.Range("J31").Value = 5905.852
$ws.Range("K31").Value = 2566.7334
$ws.Range("L31").Value = 5905.852
$ws.Range("M31").Value = -2271.7334
$ws.Range("N31").Value = -6495.852

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4713.3096
$ws.Range("I34").Value = 2566.7334
$ws.Range("J34").Value = 5905.852
$ws.Range("K34").Value = 2566.7334
$ws.Range("L34").Value = 5905.852
$ws.Range("M34").Value = -2364.7334
$ws.Range("N34").Value = -6309.852

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1110.72
$ws.Range("I122").Value = 777.44446
$ws.Range("K122").Value = 2332.33338
$ws.Range("M122").Value = 117.66662

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H137").Value = 28431.666
$ws.Range("J137").Value = 28431.666
$ws.Range("L137").Value = 28431.666
$ws.Range("N137").Value = -38631.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 27844.5
$ws.Range("J141").Value = 27844.5
$ws.Range("L141").Value = 27844.5
$ws.Range("N141").Value = -38204.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 3457.7646
$ws.Range("I94").Value = 1146.75
$ws.Range("J94").Value = 5512
$ws.Range("K94").Value = 3440.25
$ws.Range("L94").Value = 16536
$ws.Range("M94").Value = -2764.25
$ws.Range("N94").Value = -17888

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 100000500
$ws.Range("I118").Value = 166666830
$ws.Range("J118").Value = 1000
$ws.Range("K118").Value = 500000490
$ws.Range("L118").Value = 3000
$ws.Range("M118").Value = -499999247
$ws.Range("N118").Value = -5486

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 696.16
$ws.Range("J131").Value = 743.6222
$ws.Range("L131").Value = 2230.8666
$ws.Range("N131").Value = -12310.8666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 919.46155
$ws.Range("I132").Value = 826.5
$ws.Range("J132").Value = 999.1429000000001
$ws.Range("K132").Value = 7438.5
$ws.Range("L132").Value = 8992.286100000001
$ws.Range("M132").Value = -4908.5
$ws.Range("N132").Value = -14052.2861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3595.84
$ws.Range("I80").Value = 3533.3333
$ws.Range("J80").Value = 3631
$ws.Range("K80").Value = 3533.3333
$ws.Range("L80").Value = 3631
$ws.Range("M80").Value = -2535.3333
$ws.Range("N80").Value = -5627

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3595.84
$ws.Range("I83").Value = 3533.3333
$ws.Range("J83").Value = 3631
$ws.Range("K83").Value = 17666.6665
$ws.Range("L83").Value = 18155
$ws.Range("M83").Value = -12674.6665
$ws.Range("N83").Value = -28139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 33950.055
$ws.Range("I132").Value = 7691.273
$ws.Range("J132").Value = 75213.86
$ws.Range("K132").Value = 23073.819
$ws.Range("L132").Value = 225641.58
$ws.Range("M132").Value = -20543.819
$ws.Range("N132").Value = -230701.58

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 28050.666
$ws.Range("J134").Value = 28050.666
$ws.Range("L134").Value = 84151.99800000001
$ws.Range("N134").Value = -89221.99800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5283.3335
$ws.Range("I7").Value = 5309.091
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 5309.091
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -5197.091
$ws.Range("N7").Value = -5224

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4010.4707
$ws.Range("I40").Value = 3620.9285
$ws.Range("J40").Value = 5828.3335
$ws.Range("K40").Value = 3620.9285
$ws.Range("L40").Value = 5828.3335
$ws.Range("M40").Value = -3484.9285
$ws.Range("N40").Value = -6100.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2270.0667
$ws.Range("I93").Value = 2087.5833
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 2087.5833
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = -839.5832999999998
$ws.Range("N93").Value = -5496

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5283.3335
$ws.Range("I126").Value = 5309.091
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 15927.273
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -13457.273
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1299.5834
$ws.Range("I122").Value = 1010
$ws.Range("K122").Value = 3030
$ws.Range("M122").Value = -580

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 80715
$ws.Range("J141").Value = 80715
$ws.Range("L141").Value = 80715
$ws.Range("N141").Value = -91075

